# Append rows 230-233 (new daily data through 2021-04-21) to Sheet1,
# mirroring the existing row layout: column A holds the date (styled
# like the preceding rows), columns B:AX hold the per-comune counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 229

$newRows = @(
    ,@(44304,@(2,1,0,6,18,4,2,3,3,0,4,10,7,0,1,0,10,3,5,4,35,1,1,3,7,0,3,2,3,1,0,18,1,2,1,4,2,7,1,6,185,1,0,2,0,1,0,0,0))
    ,@(44305,@(3,0,1,10,6,1,0,1,1,0,0,11,6,0,0,0,5,2,1,5,59,0,0,1,4,0,0,1,0,11,0,13,2,0,0,5,4,4,0,3,161,0,0,0,1,0,0,0,0))
    ,@(44306,@(2,1,0,14,10,0,4,0,2,0,1,5,2,0,0,0,3,0,0,2,17,0,1,5,11,0,0,3,0,3,2,11,2,11,3,3,7,8,0,1,135,0,0,0,0,0,0,1,0))
    ,@(44307,@(6,0,2,2,5,0,0,0,1,0,1,6,4,0,0,0,7,0,0,0,11,0,0,1,0,2,5,0,1,0,1,15,0,0,1,2,0,4,0,5,83,0,0,0,0,0,1,0,0))
)

foreach ($rowData in $newRows) {
    $lastDataRow++
    $dateValue = $rowData[0]
    $counts = $rowData[1]

    # Carry the date column style (border/bold/center + custom date
    # number format) down from the previous row, then overwrite the value.
    $ws.Range("A" + ($lastDataRow - 1)).Copy($ws.Range("A" + $lastDataRow))
    $ws.Range("A" + $lastDataRow).Value = $dateValue

    for ($i = 0; $i -lt $counts.Length; $i++) {
        $ws.Cells.Item($lastDataRow, $i + 2).Value = $counts[$i]
    }
}

Write-Output ("Added rows 230-233; last row used is now " + $lastDataRow)